# Update "心之器等级" sheet: fill in the real level/breakthrough values
# (replacing the 5 / -1 placeholders), then sort the table descending by
# column B (Data > Sort), which is what produced the AutoFilter + sortState
# left behind in the saved workbook.

$wb = $excel.ActiveWorkbook
$wsLevel = $wb.Worksheets.Item(1)   # 心之器等级
$wsAttr  = $wb.Worksheets.Item(2)   # 心之器属性

# New "突破数" (col B) and "等级" (col C) values, in the ORIGINAL (pre-sort)
# row order -- i.e. row 2 .. row 56 of the sheet as it exists today.
$bVals = @(-1,5,0,-1,-1,1,1,-1,-1,1,0,-1,-1,2,2,5,3,2,-1,0,4,2,5,-1,5,5,-1,0,5,5,5,-1,5,-1,5,5,5,5,-1,-1,-1,-1,-1,-1,-1,-1,-1,-1,-1,-1,-1,-1,-1,-1,-1)
$cVals = @(0,40,25,0,0,30,30,0,0,30,25,0,0,35,35,40,40,35,0,25,40,35,40,0,40,40,0,25,40,40,40,0,40,0,40,40,40,40,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

for ($i = 0; $i -lt $bVals.Length; $i++) {
    $r = $i + 2
    $wsLevel.Cells.Item($r, 2).Value = $bVals[$i]
    $wsLevel.Cells.Item($r, 3).Value = $cVals[$i]
}

# Sort A1:C56 (header row included) descending on column B -- reproduces
# the AutoFilter/sortState block Excel writes when you run Data > Sort.
$sortRange = $wsLevel.Range("A1:C56")
$keyRange = $wsLevel.Range("B1:B56")
$wsLevel.Sort.SortFields.Clear()
$wsLevel.Sort.SortFields.Add($keyRange, 0, 2, 0, 0)
$wsLevel.Sort.SetRange($sortRange)
$wsLevel.Sort.Header = 1
$wsLevel.Sort.Orientation = 1
$wsLevel.Sort.Apply()

# Column widths (best-fit, as Excel leaves them after auto-fitting A:C).
$wsLevel.Columns.Item(1).ColumnWidth = 19.21875
$wsLevel.Columns.Item(2).ColumnWidth = 7.5546875
$wsLevel.Columns.Item(3).ColumnWidth = 5.5546875

# Make 心之器等级 the active/selected sheet with B4 selected, and drop the
# "tabSelected" flag from 心之器属性 (it was the active sheet before).
$wsLevel.Activate()
$wsLevel.Range("B4").Select()
